$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("98_2")

$ws.Range("C8").Value = 12291
$ws.Range("C8").NumberFormat = $ws.Range("B8").NumberFormat
$ws.Range("D11").Value = 41
$ws.Range("C24").Value = 335
